# Refined metadata to be additional tab
#
# 1. Update the "time_taken" timestamps (column F) on the "data" sheet.
# 2. Add a new "metadata" worksheet (placed after "data") describing the
#    PanelApp query that produced the "data" sheet.

$wb = $excel.ActiveWorkbook
$data = $wb.Worksheets.Item(1)

# ---------------------------------------------------------------------
# 1. Refresh the time_taken column on the "data" sheet with the new
#    query timestamps.
# ---------------------------------------------------------------------
$newTimes = @(
  "2021-10-05 14:22:23.195051",
  "2021-10-05 14:22:23.195057",
  "2021-10-05 14:22:23.195060",
  "2021-10-05 14:22:23.195061",
  "2021-10-05 14:22:23.195064",
  "2021-10-05 14:22:23.195066",
  "2021-10-05 14:22:23.195068",
  "2021-10-05 14:22:23.195070",
  "2021-10-05 14:22:23.195072",
  "2021-10-05 14:22:23.195074",
  "2021-10-05 14:22:23.195076",
  "2021-10-05 14:22:23.195077",
  "2021-10-05 14:22:23.195079",
  "2021-10-05 14:22:23.195081",
  "2021-10-05 14:22:23.195083",
  "2021-10-05 14:22:23.195085",
  "2021-10-05 14:22:23.195087",
  "2021-10-05 14:22:23.195089",
  "2021-10-05 14:22:23.195091",
  "2021-10-05 14:22:23.195093",
  "2021-10-05 14:22:23.195095",
  "2021-10-05 14:22:23.195097",
  "2021-10-05 14:22:23.195099",
  "2021-10-05 14:22:23.195101",
  "2021-10-05 14:22:23.195103",
  "2021-10-05 14:22:23.195129"
)

for ($i = 0; $i -lt $newTimes.Length; $i++) {
    $row = $i + 2
    $data.Cells.Item($row, 6).Value = $newTimes[$i]
}

# ---------------------------------------------------------------------
# 2. Add the "metadata" worksheet right after "data".
# ---------------------------------------------------------------------
$meta = $wb.Worksheets.Add($null, $data)
$meta.Name = "metadata"

$meta.Cells.Item(1, 2).Value = "data_name"
$meta.Cells.Item(1, 3).Value = "data_id"
$meta.Cells.Item(1, 4).Value = "data_version"
$meta.Cells.Item(1, 5).Value = "data_version_created"
$meta.Cells.Item(1, 6).Value = "panel_query_time"
$meta.Cells.Item(1, 7).Value = "panel_get_request"

$meta.Cells.Item(2, 1).Value = 0
$meta.Cells.Item(2, 2).Value = "Pyruvate dehydrogenase (PDH) deficiency"
$meta.Cells.Item(2, 3).Value = 531
$meta.Cells.Item(2, 4).Value = "'1.30"
$meta.Cells.Item(2, 5).Value = "2021-05-26T12:00:49.670416Z"
$meta.Cells.Item(2, 6).Value = "2021-10-05 14:22:23.192718"
$meta.Cells.Item(2, 7).Value = "https://panelapp.genomicsengland.co.uk/api/v1/panels/531/?format=json"

# Match styling of the "data" sheet: bold header row + bold leading index
# column (reuse the same cell style already present in the workbook by
# copying formats from the "data" sheet's header cell).
$data.Range("B1").Copy()
$meta.Range("B1:G1").PasteSpecial(-4122)  # xlPasteFormats
$meta.Cells.Item(2, 1).PasteSpecial(-4122)
$excel.CutCopyMode = $false

$meta.Select()
$data.Select()
